$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read all affected source rows (B:AC) first, to avoid overwriting data before it is read.
$row5 = $ws.Range("B5:AC5").Value2
$row6 = $ws.Range("B6:AC6").Value2
$row16 = $ws.Range("B16:AC16").Value2
$row17 = $ws.Range("B17:AC17").Value2
$row25 = $ws.Range("B25:AC25").Value2
$row26 = $ws.Range("B26:AC26").Value2
$row37 = $ws.Range("B37:AC37").Value2
$row38 = $ws.Range("B38:AC38").Value2
$row85 = $ws.Range("B85:AC85").Value2
$row86 = $ws.Range("B86:AC86").Value2
$row88 = $ws.Range("B88:AC88").Value2
$row89 = $ws.Range("B89:AC89").Value2
$row92 = $ws.Range("B92:AC92").Value2
$row93 = $ws.Range("B93:AC93").Value2
$row98 = $ws.Range("B98:AC98").Value2
$row99 = $ws.Range("B99:AC99").Value2
$row100 = $ws.Range("B100:AC100").Value2
$row102 = $ws.Range("B102:AC102").Value2
$row103 = $ws.Range("B103:AC103").Value2
$row104 = $ws.Range("B104:AC104").Value2
$row107 = $ws.Range("B107:AC107").Value2
$row108 = $ws.Range("B108:AC108").Value2
$row119 = $ws.Range("B119:AC119").Value2
$row120 = $ws.Range("B120:AC120").Value2

# Write back rows in their new (swapped/rotated) order.
$ws.Range("B5:AC5").Value = $row6
$ws.Range("B6:AC6").Value = $row5

$ws.Range("B16:AC16").Value = $row17
$ws.Range("B17:AC17").Value = $row16

$ws.Range("B25:AC25").Value = $row26
$ws.Range("B26:AC26").Value = $row25

$ws.Range("B37:AC37").Value = $row38
$ws.Range("B38:AC38").Value = $row37

$ws.Range("B85:AC85").Value = $row86
$ws.Range("B86:AC86").Value = $row85

$ws.Range("B88:AC88").Value = $row89
$ws.Range("B89:AC89").Value = $row88

$ws.Range("B92:AC92").Value = $row93
$ws.Range("B93:AC93").Value = $row92

$ws.Range("B98:AC98").Value = $row99
$ws.Range("B99:AC99").Value = $row100
$ws.Range("B100:AC100").Value = $row98

$ws.Range("B102:AC102").Value = $row104
$ws.Range("B103:AC103").Value = $row102
$ws.Range("B104:AC104").Value = $row103

$ws.Range("B107:AC107").Value = $row108
$ws.Range("B108:AC108").Value = $row107

$ws.Range("B119:AC119").Value = $row120
$ws.Range("B120:AC120").Value = $row119

